$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns that receive updated values, in order
$cols = @("B", "C", "D", "E", "F", "H", "I", "K", "N")

# New values for rows 2..25, one array per row, matching $cols order
$data = @(
    @(11.63507054624224, 7.930538139316953, 5.896791424071353, 16.35060966286246, 31.76627013736633, 7.344005520526261, 24.76664481412136, 12.1168156590501, 19.7676499075444),
    @(11.32922955996976, 7.622261269460802, 5.908790502108486, 15.43104943582788, 31.59217834263301, 7.344005520526261, 24.76072662657392, 11.89048108504723, 19.82133918496908),
    @(11.14079757440247, 7.429163852767115, 5.916309098691179, 14.8428957536524, 31.49368683623382, 7.344005520526261, 24.76249397089042, 11.75299025392147, 19.85620113457726),
    @(11.06397820849839, 7.349654403407249, 5.919411473706312, 14.5975735728877, 31.45569032339773, 7.344005520526261, 24.76456984480237, 11.69741583384393, 19.87088472241151),
    @(11.05122439681561, 7.336406663743994, 5.919928959359462, 14.5565061471582, 31.44951097298134, 7.344005520526261, 24.76499633333539, 11.688217580295, 19.87335173953781),
    @(11.13976151449102, 7.42809469055926, 5.916350781975117, 14.83960970666797, 31.4931657061899, 7.344005520526261, 24.76251648180601, 11.75223881016096, 19.85639723081888),
    @(11.52982676912367, 7.825114753154784, 5.900897637999581, 16.03858008405627, 31.70451677715083, 7.344005520526261, 24.76348193371844, 12.03851762694499, 19.78576844610426),
    @(12.2842034851214, 8.568011937093429, 5.871770276236498, 18.23143819052415, 32.18429079094876, 7.344005520526261, 24.80831289210481, 12.6080691574256, 19.66230567844637),
    @(12.82523132443718, 9.08578750650805, 5.851054636582391, 19.85965522486208, 32.57457172557516, 7.344005520526261, 24.86749567343276, 13.02680132400799, 19.58075486323613),
    @(13.06720586320496, 9.314169961490324, 5.841771971579065, 20.5589751465383, 32.75983199007929, 7.344005520526261, 24.90011683907963, 13.21640717235754, 19.54564177087198),
    @(13.15814373377797, 9.399548520056509, 5.838276587964648, 20.81787895845547, 32.83104958979783, 7.344005520526261, 24.91328740250825, 13.28800505306842, 19.53263075300217),
    @(13.13859089982836, 9.381210961838478, 5.839028510723185, 20.76238184345388, 32.81566506756827, 7.344005520526261, 24.91041455932931, 13.2725953040495, 19.53542020608925),
    @(13.07470178484452, 9.321216649589434, 5.841484011196629, 20.58039363580265, 32.76567006392607, 7.344005520526261, 24.90118401603781, 13.22230203576066, 19.54456562288518),
    @(13.03547491474039, 9.284322430181295, 5.842990634199905, 20.46815152124613, 32.73518376861073, 7.344005520526261, 24.89563645370004, 13.19146750542508, 19.550204646202),
    @(12.80932589292043, 9.070711453911368, 5.851664073733141, 19.81312363710412, 32.56261594584297, 7.344005520526261, 24.86547833813966, 13.01438593577, 19.58308953642506),
    @(12.66945833421393, 8.937777874002832, 5.857020693160666, 19.40071560867562, 32.45869598927803, 7.344005520526261, 24.84843583829754, 12.90547246716318, 19.60377179033026),
    @(12.58862565573052, 8.860648260627467, 5.860114979892858, 19.15961396647396, 32.39965351607954, 7.344005520526261, 24.83916998708542, 12.84275002351646, 19.61585448263216),
    @(12.56119414316173, 8.834420952439856, 5.861164952218536, 19.07731098687171, 32.37978942634479, 7.344005520526261, 24.83612492252372, 12.82150219608564, 19.61997755514871),
    @(12.68438800407232, 8.951998777527038, 5.856449099054286, 19.44502005921335, 32.46968327158296, 7.344005520526261, 24.85019452200045, 12.91707515436254, 19.60155079453317),
    @(13.09348709746871, 9.33886897486607, 5.840762239003787, 20.6340081984049, 32.78032632544519, 7.344005520526261, 24.90387307861332, 13.23708044656716, 19.54187164131724),
    @(13.35678052385573, 9.585242468286809, 5.830624905838371, 21.37662842324282, 32.98952565969058, 7.344005520526261, 24.94371948928167, 13.44501415903962, 19.50453234551685),
    @(13.21665882075062, 9.454362652382155, 5.836025045048142, 20.98341901854132, 32.87732284845554, 7.344005520526261, 24.92201758795334, 13.33417056200213, 19.52430867593997),
    @(12.67763960402961, 8.945571695017684, 5.856707470905453, 19.42500248796415, 32.46471372946995, 7.344005520526261, 24.84939776343443, 12.91182990582755, 19.60255430804276),
    @(12.08198253152986, 8.371555263831036, 5.879527460027592, 17.63262602567063, 32.04770359443593, 7.344005520526261, 24.79157776242099, 12.45360935518256, 19.69409679287172)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowValues[$j]
    }
}
